$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 630 (pushes old rows 630-671 down to 633-674)
$ws.Range("A630:A632").EntireRow.Insert()

$data = @(
  @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44706, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Pintón", 80, 10000, 10000, 10000, "`$/caja 20 kilos", "Ecuador", 500, 20),
  @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44706, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Primera Maduro", 120, 12000, 12000, 12000, "`$/caja 20 kilos", "Ecuador", 600, 20),
  @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44706, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Primera Pintón", 120, 13000, 13000, 13000, "`$/caja 20 kilos", "Ecuador", 650, 20)
)

for ($i = 0; $i -lt 3; $i++) {
  $row = 630 + $i
  $vals = $data[$i]
  for ($c = 0; $c -lt 20; $c++) {
    $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
  }
}

"Done"
